$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Delete the bold run containing "12/23/2020" at the end of the
#    "Final Approval Date" line (the trailing <w:tab/> run is kept).
# ------------------------------------------------------------------
$dateRng = $d.Content
$dateRng.Find.Execute("12/23/2020", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$dateRng.Delete()

# ------------------------------------------------------------------
# 2. Remove the existing (misplaced) "_GoBack" bookmark that currently
#    sits at the end of the "... status:  " paragraph.
# ------------------------------------------------------------------
$d.Bookmarks("_GoBack").Delete()

# ------------------------------------------------------------------
# 3. Re-insert the "_GoBack" bookmark (empty span) right after the
#    "  " run that ends the "1. Title of proposed study:" paragraph,
#    without splitting/altering that existing run.
#
#    Trick: temporarily append a throw-away character after the run,
#    wrap a bookmark around just that character, then clear the
#    character's text. The empty run left behind is pruned
#    automatically, leaving a clean, adjacent bookmarkStart/bookmarkEnd
#    pair immediately after the untouched "  " run.
# ------------------------------------------------------------------
$titleRng = $d.Content
$titleRng.Find.Execute("1. Title of proposed study:  ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$titleRng.Collapse(0)
$titleRng.InsertAfter("X")
$d.Bookmarks.Add("_GoBack", $titleRng)
$titleRng.Text = ""
